$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update simple Property/Value pairs in place ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# --- Update the Contact rows with the new display text ---
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the Contact rows, pushing
#     Description/Purpose/Copyright/Immutable down by one row ---
$ws.Range("A12:B15").Copy($ws.Range("A13:B16"))
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
